# Refactored: sliding-window IPC PO results for this run, now written
# directly from the results table produced by the (relocated) analysis
# script instead of being recomputed cell-by-cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-point results: Row, IPC PO (C), DELTA = C - B (D), DELTA^2 (E)
$newResults = @(
    @{ Row = 2; C = 27.49241226758658; D = -1.48758773241342; E = 2.212917261626901 }
    @{ Row = 3; C = 28.1808189083947; D = -0.9691810916053001; E = 0.9393119883252411 }
    @{ Row = 4; C = 29.85633800906434; D = 0.5063380090643363; E = 0.2563781794232359 }
    @{ Row = 5; C = 29.60282396639035; D = 0.2328239663903524; E = 0.05420699932573595 }
    @{ Row = 6; C = 29.15167570878408; D = -0.3883242912159233; E = 0.1507957551483492 }
    @{ Row = 7; C = 28.50871812765812; D = -1.04128187234188; E = 1.084267937667812 }
    @{ Row = 8; C = 29.53306944917504; D = -0.2169305508249622; E = 0.04705886388122149 }
    @{ Row = 9; C = 30.23630653117635; D = 0.3963065311763465; E = 0.1570588666530285 }
    @{ Row = 10; C = 29.98318372129723; D = 0.1731837212972316; E = 0.02999260132235718 }
    @{ Row = 11; C = 31.1543561372367; D = 1.234356137236695; E = 1.523635073533894 }
    @{ Row = 12; C = 29.96703474200553; D = -0.01296525799446968; E = 0.00016809791486316 }
    @{ Row = 13; C = 30.03303552117937; D = -0.006964478820627562; E = 0.00004850396524296988 }
    @{ Row = 14; C = 30.09177113028868; D = -0.1182288697113236; E = 0.01397806563321714 }
    @{ Row = 15; C = 29.59311683979067; D = -0.6268831602093314; E = 0.3929824965540383 }
    @{ Row = 16; C = 30.44348981295915; D = 0.06348981295915124; E = 0.004030956349588008 }
    @{ Row = 17; C = 30.48743912317587; D = 0.04743912317586663; E = 0.002250470407695046 }
    @{ Row = 18; C = 30.45050122073281; D = -0.02949877926718969; E = 0.0008701779782543803 }
    @{ Row = 19; C = 31.57859063932618; D = 0.8885906393261784; E = 0.7895933242981066 }
    @{ Row = 20; C = 30.82411957430374; D = 0.07411957430373661; E = 0.005493711294967132 }
    @{ Row = 21; C = 32.01803992955325; D = 1.078039929553253; E = 1.162170089711182 }
    @{ Row = 22; C = 31.06554974563931; D = 0.11554974563931; E = 0.01335174371730923 }
    @{ Row = 23; C = 31.68144431219928; D = 0.6614443121992757; E = 0.4375085781407729 }
    @{ Row = 24; C = 31.93836282950241; D = 0.8183628295024086; E = 0.6697177207111883 }
    @{ Row = 25; C = 32.68288845627483; D = 1.40288845627483; E = 1.968096020749176 }
    @{ Row = 26; C = 31.92775612296992; D = 0.5477561229699184; E = 0.3000367702510364 }
    @{ Row = 27; C = 32.50689535741705; D = 0.9268953574170524; E = 0.8591350036012854 }
    @{ Row = 28; C = 31.6804080850526; D = 0.03040808505260628; E = 0.0009246516365665376 }
    @{ Row = 29; C = 32.94701976420416; D = 1.067019764204158; E = 1.138531177202297 }
    @{ Row = 30; C = 33.16377948120398; D = 0.8837794812039803; E = 0.7810661713971766 }
    @{ Row = 31; C = 32.28287875391908; D = -0.1671212460809244; E = 0.02792951089164088 }
    @{ Row = 32; C = 34.2162396324386; D = 1.366239632438599; E = 1.866610733245959 }
    @{ Row = 33; C = 32.36611826104788; D = -0.5338817389521182; E = 0.2850297111865377 }
    @{ Row = 34; C = 32.56279035552895; D = -0.5372096444710479; E = 0.2885942021127097 }
    @{ Row = 35; C = 33.30642365044577; D = -0.09357634955423322; E = 0.008756533195896043 }
    @{ Row = 36; C = 33.70611770471226; D = 0.006117704712252703; E = 0.00003742631094631893 }
    @{ Row = 37; C = 34.61605223088247; D = 0.5160522308824653; E = 0.2663099049987693 }
    @{ Row = 38; C = 34.31129008754729; D = -0.08870991245270687; E = 0.007869448567366918 }
    @{ Row = 39; C = 35.30499002121963; D = 0.404990021219632; E = 0.164016917287478 }
    @{ Row = 40; C = 35.10723594525111; D = -0.1927640547488849; E = 0.0371579808032311 }
    @{ Row = 41; C = 35.6412975690336; D = -0.05870243096639882; E = 0.003445975401364819 }
    @{ Row = 42; C = 36.21455944135352; D = -0.08544055864648215; E = 0.007300089061822955 }
    @{ Row = 43; C = 36.49702268034099; D = -0.3029773196590071; E = 0.09179525622775617 }
    @{ Row = 44; C = 37.20155162870451; D = -0.09844837129548267; E = 0.009692081810733217 }
    @{ Row = 45; C = 39.15502187870906; D = 1.255021878709059; E = 1.575079916038415 }
    @{ Row = 46; C = 39.24086098913765; D = 0.7408609891376514; E = 0.5488750052260193 }
    @{ Row = 47; C = 39.07086653964407; D = 0.1708665396440736; E = 0.02919537436993978 }
    @{ Row = 48; C = 39.51592348560242; D = 0.1159234856024227; E = 0.01343825451421511 }
    @{ Row = 49; C = 40.580271733598; D = 0.6802717335979978; E = 0.4627696315324252 }
    @{ Row = 50; C = 36.69296421120788; D = -3.407035788792122; E = 11.60789286611036 }
    @{ Row = 51; C = 39.81737960129512; D = -0.7826203987048785; E = 0.612494688468983 }
)

foreach ($r in $newResults) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# Summary rows: TOTAL delta (C52), sum of squared deltas (E52), MSE (E53)
$ws.Range("C52").Value = 5.158801916162126
$ws.Range("E52").Value = 32.9098687657843
$ws.Range("E53").Value = 0.658197375315686
